# Zeitplan.xlsx - Abgabetermin ("Definitive Abgabe der Maturaarbeit") im
# Zeitplan ergaenzt: eine neue Zeile wird vor der bisherigen "?"-Zeile
# (Schlusspraesentation) eingefuegt, mit dem Datum der darauffolgenden Woche
# und dem neuen Meilenstein-Text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Eine neue Zeile bei 18 einfuegen - alles ab der bisherigen Zeile 18
# (die "?"/Schlusspraesentation-Zeile und alles danach) rutscht eine Zeile
# weiter nach unten.
$ws.Rows("18:18").Insert()

# Neues Datum (2011-12-19, eine Woche nach dem letzten Termin 40889) und der
# neue Meilenstein-Text in der frisch eingefuegten Zeile.
$ws.Range("A18").Value = 40896
$ws.Range("A18").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"
$ws.Range("B18").Value = "Definitive Abgabe der Maturaarbeit"

# Auswahl im Blatt auf die neu gepflegte Zelle verschieben.
$ws.Range("B18").Select()
